$d = $word.ActiveDocument

function Find-ParagraphContaining {
    param([string]$NeedleText)

    $paras = $d.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($NeedleText)) {
            return $p.Range
        }
    }
    return $null
}

function Replace-ParagraphRun {
    # Replaces the run content of a whole paragraph (everything up to, but not
    # including, the trailing paragraph mark) with a fresh sequence of <w:r>
    # elements, so the pieces land as distinct, unmerged runs in the OOXML.
    param([string]$NeedleText, [string]$RunsXml)

    $paraRange = Find-ParagraphContaining $NeedleText
    $target = $d.Range($paraRange.Start, $paraRange.End - 1)

    $xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$RunsXml</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

    $target.InsertXML($xml)
}

function Build-TextRun {
    param([string]$RunPropsXml, [string]$Text)
    $escaped = $Text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $needsPreserve = ($Text.Length -eq 0) -or ($Text -ne $Text.Trim())
    if ($needsPreserve) {
        return "<w:r><w:rPr>$RunPropsXml</w:rPr><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    } else {
        return "<w:r><w:rPr>$RunPropsXml</w:rPr><w:t>$escaped</w:t></w:r>"
    }
}

$smallProps = '<w:rFonts w:cs="TH SarabunPSK" w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK"/><w:sz w:val="10"/><w:szCs w:val="10"/>'
$bigProps   = '<w:rFonts w:cs="TH SarabunPSK" w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/>'

# {#activities}  ->  {#  /  activity_image_captions  /  }
$runs1 = (Build-TextRun $smallProps "{#") + (Build-TextRun $smallProps "activity_image_captions") + (Build-TextRun $smallProps "}")
Replace-ParagraphRun "{#activities}" $runs1

# <tab/><tab/>{order}. {activity_name_govt}  ->  <tab/><tab/>{order}. {  /  caption  /  }
$runs2 = "<w:r><w:rPr>$bigProps</w:rPr><w:tab/><w:tab/></w:r>" + (Build-TextRun $bigProps "{order}. {") + (Build-TextRun $bigProps "caption") + (Build-TextRun $bigProps "}")
Replace-ParagraphRun "{order}. {activity_name_govt}" $runs2

# {/activities}  ->  {/  /  activity_image_captions  /  }
$runs3 = (Build-TextRun $smallProps "{/") + (Build-TextRun $smallProps "activity_image_captions") + (Build-TextRun $smallProps "}")
Replace-ParagraphRun "{/activities}" $runs3
